$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update day 1 (July 2025) total_venda value
$ws.Range("B2").Value = 17734.16

# Insert a new row for day 2 (July 2025) right after row 2, shifting
# everything below down by one row.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20464.65
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 2025
$ws.Range("E3").Value = "07/2025"
